# Update the "Total" row of the marksheet: correct answers went from 3 to 5
# (Marking row) and the resulting raw score / max went from 69/84 to
# 115/140 (Total row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 115
$ws.Range("E12").Value = "115/140"
